# Update the existing "scheme name" labels in column B for rows 10-16.
# (Rows 3-9 keep their original labels.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# Append three new rows (17-19) of averaged-intensity results, matching the
# existing table's layout: column A holds the numeric row index (styled like
# the other index cells), column B holds the scheme name, and columns C:M
# hold the per-HKL averaged intensity values (all 1 for this data set).

$newRows = @(
    @{ Row = 17; Index = 15; Name = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Index = 16; Name = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Index = 17; Name = "HexGrid-60degTilt5degRes" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the formatting from the row above (A16 has the bold/bordered
    # "index" style) so the new index cell matches the rest of the column.
    $ws.Range("A16").Copy() | Out-Null
    $ws.Range("A$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$rowNum").Value = $r.Index
    $ws.Range("B$rowNum").Value = $r.Name

    $ws.Range("C$rowNum").Value = 1
    $ws.Range("D$rowNum").Value = 1
    $ws.Range("E$rowNum").Value = 1
    $ws.Range("F$rowNum").Value = 1
    $ws.Range("G$rowNum").Value = 1
    $ws.Range("H$rowNum").Value = 1
    $ws.Range("I$rowNum").Value = 1
    $ws.Range("J$rowNum").Value = 1
    $ws.Range("K$rowNum").Value = 1
    $ws.Range("L$rowNum").Value = 1
    $ws.Range("M$rowNum").Value = 1
}

$excel.CutCopyMode = 0
